# Rename the "*img" sheets to "img*" (e.g. "himg" -> "imgh", "eimg" -> "imge")
# and activate the last sheet ("imge", formerly "eimg"), matching the diff's
# activeTab/tabSelected change (active tab moves from index 3 "xbday" to
# index 16 "imge").

$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# The workbook's active tab moves to the last sheet ("imge", formerly "eimg"),
# which also sets that sheet's tabSelected=true and clears it on the
# previously active sheet ("xbday").
$wb.Worksheets.Item("imge").Activate()
